$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 6: reuse the existing "full border" row style (like rows 2-5) ---
$ws.Range("A2:D2").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Rows.Item(6).RowHeight = 40.2

# --- New rows 7 & 8: style copied from an existing B-cell, then edges trimmed ---
$ws.Range("B4").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Borders.Item(8).LineStyle = -4142
$ws.Range("B7").Borders.Item(9).LineStyle = -4142
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Rows.Item(7).RowHeight = 39.6
$ws.Rows.Item(8).RowHeight = 39.6

# --- Fill in the new/updated question text (order matters for shared-string layout) ---
$ws.Range("B6").Value = "Should admin be able to delete/add/update users and comments?"
$ws.Range("B7").Value = "Should rating system appears as 5 stars rating in the history"
$ws.Range("B8").Value = "Data will be collected and sent  to third party to check the availability?"

$ws.Range("B2").Value = "Should registration form contains user name/ email /phone number?"
$ws.Range("B3").Value = "Should username contains special character,at least 3 letters?"
$ws.Range("B5").Value = "Should the password have at least 8 characters and at least 1 special character?"
$ws.Range("B4").Value = "Should the user log in with email and password?"

# Row 4's wording got shorter -> its wrapped height shrinks
$ws.Rows.Item(4).RowHeight = 27

# --- Data validations ---
$cUnion = $excel.Union($ws.Range("C4"), $ws.Range("C5"), $ws.Range("C6"))
$cUnion.Validation.Add(3, 1, 1, '"Yes , No"')
$ws.Range("C2:C3").Validation.Add(3, 1, 1, '"Yes, No"')

# --- Print/page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to the next empty row ---
$ws.Range("B9").Select()
